# Auto-generated edit script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '35.339.68'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.883.06'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '245.42'
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +2.78%  '
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '54.95'
$ws.Range('E10').Value = '  +7.14%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0742'
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0979'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '13.63'
$ws.Range('E13').Value = '  +5.95%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.158.20'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('E15').Value = '  +7.19%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.99'
$ws.Range('E16').Value = '  +1.56%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.885.39'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '35.310.63'
$ws.Range('E18').Value = '  +1.55%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '73.23'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0824'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '243.75'
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '12.77'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('E23').Value = '  +4.27%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.62'
$ws.Range('E24').Value = '  +8.92%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  -4.71%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.28'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.57'
$ws.Range('E28').Value = '  +2.23%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.22'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('E31').Value = '  +2.90%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.27'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  +13.73%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.14'
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  -13.31%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.846'
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('E38').Value = '  -2.88%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0719'
$ws.Range('E39').Value = '  +8.93%  '
$ws.Range('E40').Value = '  +3.53%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '97.35'
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '17.06'
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.53'
$ws.Range('E44').Value = '  +10.90%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.310.15'
$ws.Range('E45').Value = '  +2.01%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.37'
$ws.Range('E46').Value = '  +1.23%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0808'
$ws.Range('E47').Value = '  +3.54%  '
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.28'
$ws.Range('E50').Value = '  -2.98%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.057.58'
$ws.Range('E51').Value = '  -0.05%  '
